# Unit Test Protocol (Moving Car) - add missing test-step detail for rows 4-9
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Backlog (Moving Car)")

# Row 4: Timer 0 Driver test case
$ws.Range("D4").Value = "Call TIMER0_delay(3000) to delay the code for 3 seconds"
$ws.Range("E4").Value = "delay 3 seconds"
$ws.Range("F4").Value = "code is delayed for 3 seconds"

# Row 5: Timer 2 Driver test case
$ws.Range("D5").Value = "Call TIMER0_delay(3000) to delay the code for 3 seconds"
$ws.Range("E5").Value = "delay 3 seconds"
$ws.Range("F5").Value = "code is delayed for 3 seconds"

# Row 8: Motor Driver test case
$ws.Range("D8").Value = "Call MOTOR_control(enable pin,port number,duty cycle) "
$ws.Range("E8").Value = "car should move at 50% speed"
$ws.Range("F8").Value = "car moves at 50% speed"

# Row 9: Application (car_cycle() function) test case
$ws.Range("D9").Value = "call car_cycle() the main application"
$ws.Range("E9").Value = "car keeps looping over the states"
$ws.Range("F9").Value = "car kept looping over the states"

# Row 6: Button Driver test case
$ws.Range("D6").Value = "call BUTTON_read(pin number, port number) to turn on led"
$ws.Range("E6").Value = "led should turn on"
$ws.Range("F6").Value = "led turned on"

# Row 7: LED Driver test case
$ws.Range("D7").Value = "call LED_toggle(pin number, port number)"
$ws.Range("E7").Value = "led should toggle"
$ws.Range("F7").Value = "led keep toggling"

# Update the active selection to reflect where editing ended
$ws.Range("D15").Select()
